# "updated plan and code"
# - Add a note "Mit Mario besprechen" next to the abgabedatum row (C11)
# - Duplicate the "David" note down onto row 20 (C20)
# - Leave the selection on the last-edited cell (C20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$ws.Range("C11").Value = "Mit Mario besprechen"
$ws.Range("C20").Value = "David"

$ws.Range("C20").Select()
